$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells are treated as text so values like "1.00" or "388.80"
# are not auto-converted to numbers by Excel, matching the original inline-string data.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "63.401.28"
$ws.Range("E2").Value = "  -6.85%  "
$ws.Range("D3").Value = "3.534.07"
$ws.Range("E3").Value = "  -2.25%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "388.80"
$ws.Range("E5").Value = "  -6.63%  "
$ws.Range("D6").Value = "121.50"
$ws.Range("E6").Value = "  -6.39%  "
$ws.Range("D7").Value = "3.525.80"
$ws.Range("E7").Value = "  -2.08%  "
$ws.Range("D8").Value = "0.581"
$ws.Range("E8").Value = "  -11.36%  "
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  +0.09%  "
$ws.Range("D10").Value = "0.674"
$ws.Range("E10").Value = "  -11.06%  "
$ws.Range("E11").Value = "  -22.07%  "
$ws.Range("E12").Value = "  -23.12%  "
$ws.Range("D13").Value = "38.31"
$ws.Range("E13").Value = "  -8.52%  "
$ws.Range("D14").Value = "4.087.05"
$ws.Range("E14").Value = "  -2.40%  "
$ws.Range("D15").Value = "9.06"
$ws.Range("E15").Value = "  -7.37%  "
$ws.Range("E16").Value = "  -3.00%  "
$ws.Range("D17").Value = "3.524.41"
$ws.Range("E17").Value = "  -2.62%  "
$ws.Range("D18").Value = "12.65"
$ws.Range("E18").Value = "  +3.01%  "
$ws.Range("E19").Value = "  -7.07%  "
$ws.Range("D20").Value = "63.353.39"
$ws.Range("E20").Value = "  -6.73%  "
$ws.Range("D21").Value = "1.00"
$ws.Range("E21").Value = "  -9.59%  "
$ws.Range("D22").Value = "388.28"
$ws.Range("E22").Value = "  -15.31%  "
$ws.Range("D23").Value = "13.82"
$ws.Range("E23").Value = "  +4.71%  "
$ws.Range("D24").Value = "80.37"
$ws.Range("E24").Value = "  -9.47%  "
$ws.Range("D25").Value = "2.85"
$ws.Range("E25").Value = "  -6.39%  "
$ws.Range("E26").Value = "  +10.37%  "
$ws.Range("D27").Value = "33.46"
$ws.Range("E27").Value = "  -5.27%  "
$ws.Range("D28").Value = "2.95"
$ws.Range("E28").Value = "  -9.64%  "
$ws.Range("E29").Value = "  -14.43%  "
$ws.Range("B30").Value = "Cosmos"
$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D30").Value = "11.68"
$ws.Range("E30").Value = "  -3.99%  "
$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").Value = "2.60"
$ws.Range("E31").Value = "  -6.57%  "
$ws.Range("E32").Value = "  -7.19%  "
$ws.Range("D33").Value = "6.66"
$ws.Range("E33").Value = "  -7.18%  "
$ws.Range("E34").Value = "  -5.31%  "
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.24%  "
$ws.Range("D36").Value = "36.20"
$ws.Range("E36").Value = "  -9.92%  "
$ws.Range("D37").Value = "53.40"
$ws.Range("E37").Value = "  -4.39%  "
$ws.Range("D38").Value = "0.0432"
$ws.Range("E38").Value = "  -10.23%  "
$ws.Range("D39").Value = "0.997"
$ws.Range("E39").Value = "  -0.28%  "
$ws.Range("E40").Value = "  +2.70%  "
$ws.Range("D41").Value = "0.129"
$ws.Range("E41").Value = "  -11.97%  "
$ws.Range("B42").Value = "ApeXProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D42").Value = "3.03"
$ws.Range("E42").Value = "  +15.50%  "
$ws.Range("B43").Value = "Monero"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D43").Value = "139.78"
$ws.Range("E43").Value = "  -5.69%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "25.56"
$ws.Range("E44").Value = "  +20.10%  "
$ws.Range("D45").Value = "0.0₃0599"
$ws.Range("E45").Value = "  -23.44%  "
$ws.Range("D46").Value = "1.94"
$ws.Range("E46").Value = "  +0.95%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").Value = "4.04"
$ws.Range("E47").Value = "  -4.27%  "
$ws.Range("E48").Value = "  -5.93%  "
$ws.Range("B49").Value = "WEMIXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D49").Value = "2.45"
$ws.Range("E49").Value = "  -9.75%  "
$ws.Range("E50").Value = "  -10.74%  "
$ws.Range("D51").Value = "0.271"
$ws.Range("E51").Value = "  -9.66%  "
